$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 45.71598933333333
$ws.Range("H2").Value = 137.147968
$ws.Range("I2").Value = 0.6549002937372808
$ws.Range("J2").Value = 0.6549002937372808
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 47.991936
$ws.Range("N2").Value = 143.975808
$ws.Range("O2").Value = 0.7605119179168339
$ws.Range("P2").Value = 0.7605119179168338
$ws.Range("Q2").Value = 2193.998834262016
$ws.Range("R2").Value = 19745.98950835814
$ws.Range("S2").Value = 0.4980594784344373
$ws.Range("T2").Value = 0.4980594784344372
# Row 3
$ws.Range("G3").Value = 45.71598933333333
$ws.Range("H3").Value = 137.147968
$ws.Range("I3").Value = 0.6549002937372808
$ws.Range("J3").Value = 0.6549002937372808
$ws.Range("O3").Value = 0.1317597634642934
$ws.Range("P3").Value = 0.1317597634642934
$ws.Range("Q3").Value = 380.1133955075128
$ws.Range("R3").Value = 3421.020559567616
$ws.Range("S3").Value = 0.08628950779552037
$ws.Range("T3").Value = 0.08628950779552037
# Row 4
$ws.Range("G4").Value = 45.71598933333333
$ws.Range("H4").Value = 137.147968
$ws.Range("I4").Value = 0.6549002937372808
$ws.Range("J4").Value = 0.6549002937372808
$ws.Range("M4").Value = 6.744108333333334
$ws.Range("N4").Value = 20.232325
$ws.Range("O4").Value = 0.1068715953284784
$ws.Range("P4").Value = 0.1068715953284784
$ws.Range("Q4").Value = 308.3135846295111
$ws.Range("R4").Value = 2774.8222616656
$ws.Range("S4").Value = 0.06999023917279233
$ws.Range("T4").Value = 0.06999023917279232
# Row 5
$ws.Range("G5").Value = 45.71598933333333
$ws.Range("H5").Value = 137.147968
$ws.Range("I5").Value = 0.6549002937372808
$ws.Range("J5").Value = 0.6549002937372808
$ws.Range("M5").Value = 0.05406333333333333
$ws.Range("N5").Value = 0.16219
$ws.Range("O5").Value = 0.0008567232903942534
$ws.Range("P5").Value = 0.0008567232903942534
$ws.Range("Q5").Value = 2.471558769991111
$ws.Range("R5").Value = 22.24402892992
$ws.Range("S5").Value = 0.0005610683345307662
$ws.Range("T5").Value = 0.0005610683345307662
# Row 6
$ws.Range("I6").Value = 0.1818108415648851
$ws.Range("J6").Value = 0.1818108415648851
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 47.991936
$ws.Range("N6").Value = 143.975808
$ws.Range("O6").Value = 0.7605119179168339
$ws.Range("P6").Value = 0.7605119179168338
$ws.Range("Q6").Value = 609.089319800448
$ws.Range("R6").Value = 5481.803878204031
$ws.Range("S6").Value = 0.1382693118165844
$ws.Range("T6").Value = 0.1382693118165843
# Row 7
$ws.Range("I7").Value = 0.1818108415648851
$ws.Range("J7").Value = 0.1818108415648851
$ws.Range("O7").Value = 0.1317597634642934
$ws.Range("P7").Value = 0.1317597634642934
$ws.Range("S7").Value = 0.02395535347983338
$ws.Range("T7").Value = 0.02395535347983338
# Row 8
$ws.Range("I8").Value = 0.1818108415648851
$ws.Range("J8").Value = 0.1818108415648851
$ws.Range("M8").Value = 6.744108333333334
$ws.Range("N8").Value = 20.232325
$ws.Range("O8").Value = 0.1068715953284784
$ws.Range("P8").Value = 0.1068715953284784
$ws.Range("Q8").Value = 85.59280370374167
$ws.Range("R8").Value = 770.3352333336751
$ws.Range("S8").Value = 0.01943041468605251
$ws.Range("T8").Value = 0.0194304146860525
# Row 9
$ws.Range("I9").Value = 0.1818108415648851
$ws.Range("J9").Value = 0.1818108415648851
$ws.Range("M9").Value = 0.05406333333333333
$ws.Range("N9").Value = 0.16219
$ws.Range("O9").Value = 0.0008567232903942534
$ws.Range("P9").Value = 0.0008567232903942534
$ws.Range("Q9").Value = 0.6861444165566666
$ws.Range("R9").Value = 6.17529974901
$ws.Range("S9").Value = 0.0001557615824148167
$ws.Range("T9").Value = 0.0001557615824148166
# Row 10
$ws.Range("G10").Value = 11.24784666666667
$ws.Range("H10").Value = 33.74354
$ws.Range("I10").Value = 0.161130015850732
$ws.Range("J10").Value = 0.161130015850732
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 47.991936
$ws.Range("N10").Value = 143.975808
$ws.Range("O10").Value = 0.7605119179168339
$ws.Range("P10").Value = 0.7605119179168338
$ws.Range("Q10").Value = 539.80593736448
$ws.Range("R10").Value = 4858.25343628032
$ws.Range("S10").Value = 0.12254129738861
$ws.Range("T10").Value = 0.12254129738861
# Row 11
$ws.Range("G11").Value = 11.24784666666667
$ws.Range("H11").Value = 33.74354
$ws.Range("I11").Value = 0.161130015850732
$ws.Range("J11").Value = 0.161130015850732
$ws.Range("O11").Value = 0.1317597634642934
$ws.Range("P11").Value = 0.1317597634642934
$ws.Range("Q11").Value = 93.52214074249777
$ws.Range("R11").Value = 841.6992666824799
$ws.Range("S11").Value = 0.02123045277549029
$ws.Range("T11").Value = 0.02123045277549029
# Row 12
$ws.Range("G12").Value = 11.24784666666667
$ws.Range("H12").Value = 33.74354
$ws.Range("I12").Value = 0.161130015850732
$ws.Range("J12").Value = 0.161130015850732
$ws.Range("M12").Value = 6.744108333333334
$ws.Range("N12").Value = 20.232325
$ws.Range("O12").Value = 0.1068715953284784
$ws.Range("P12").Value = 0.1068715953284784
$ws.Range("Q12").Value = 75.85669643672223
$ws.Range("R12").Value = 682.7102679305
$ws.Range("S12").Value = 0.01722022184927074
$ws.Range("T12").Value = 0.01722022184927074
# Row 13
$ws.Range("G13").Value = 11.24784666666667
$ws.Range("H13").Value = 33.74354
$ws.Range("I13").Value = 0.161130015850732
$ws.Range("J13").Value = 0.161130015850732
$ws.Range("M13").Value = 0.05406333333333333
$ws.Range("N13").Value = 0.16219
$ws.Range("O13").Value = 0.0008567232903942534
$ws.Range("P13").Value = 0.0008567232903942534
$ws.Range("Q13").Value = 0.6080960836222221
$ws.Range("R13").Value = 5.4728647526
$ws.Range("S13").Value = 0.0001380438373609173
$ws.Range("T13").Value = 0.0001380438373609173
# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1507006666666667
$ws.Range("H14").Value = 0.452102
$ws.Range("I14").Value = 0.00215884884710222
$ws.Range("J14").Value = 0.00215884884710222
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 47.991936
$ws.Range("N14").Value = 143.975808
$ws.Range("O14").Value = 0.7605119179168339
$ws.Range("P14").Value = 0.7605119179168338
$ws.Range("Q14").Value = 7.232416749824001
$ws.Range("R14").Value = 65.091750748416
$ws.Range("S14").Value = 0.001641830277202255
$ws.Range("T14").Value = 0.001641830277202255
# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1507006666666667
$ws.Range("H15").Value = 0.452102
$ws.Range("I15").Value = 0.00215884884710222
$ws.Range("J15").Value = 0.00215884884710222
$ws.Range("O15").Value = 0.1317597634642934
$ws.Range("P15").Value = 0.1317597634642934
$ws.Range("Q15").Value = 1.253026412580444
$ws.Range("R15").Value = 11.277237713224
$ws.Range("S15").Value = 0.000284449413449351
$ws.Range("T15").Value = 0.0002844494134493509
# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1507006666666667
$ws.Range("H16").Value = 0.452102
$ws.Range("I16").Value = 0.00215884884710222
$ws.Range("J16").Value = 0.00215884884710222
$ws.Range("M16").Value = 6.744108333333334
$ws.Range("N16").Value = 20.232325
$ws.Range("O16").Value = 0.1068715953284784
$ws.Range("P16").Value = 0.1068715953284784
$ws.Range("Q16").Value = 1.016341621905556
$ws.Range("R16").Value = 9.147074597150002
$ws.Range("S16").Value = 0.0002307196203628607
$ws.Range("T16").Value = 0.0002307196203628606
# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1507006666666667
$ws.Range("H17").Value = 0.452102
$ws.Range("I17").Value = 0.00215884884710222
$ws.Range("J17").Value = 0.00215884884710222
$ws.Range("M17").Value = 0.05406333333333333
$ws.Range("N17").Value = 0.16219
$ws.Range("O17").Value = 0.0008567232903942534
$ws.Range("P17").Value = 0.0008567232903942534
$ws.Range("Q17").Value = 0.008147380375555556
$ws.Range("R17").Value = 0.07332642338000001
$ws.Range("S17").Value = 0.000001849536087753255
$ws.Range("T17").Value = 0.000001849536087753254
